# Auto-generated Excel COM-interop script to apply data update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 598408
$ws.Cells.Item(2, 4).Value = 154088
$ws.Cells.Item(2, 5).Value = 1024211681
$ws.Cells.Item(8, 3).Value = 2718
$ws.Cells.Item(8, 5).Value = 15895756
$ws.Cells.Item(10, 3).Value = 250166
$ws.Cells.Item(10, 4).Value = 63067
$ws.Cells.Item(10, 5).Value = 1048883391
$ws.Cells.Item(13, 3).Value = 133624
$ws.Cells.Item(13, 4).Value = 32270
$ws.Cells.Item(13, 5).Value = 597432994
$ws.Cells.Item(16, 3).Value = 7499
$ws.Cells.Item(16, 5).Value = 17193543
$ws.Cells.Item(19, 3).Value = 18186
$ws.Cells.Item(19, 4).Value = 3997
$ws.Cells.Item(19, 5).Value = 64613237
$ws.Cells.Item(21, 3).Value = 138805
$ws.Cells.Item(21, 4).Value = 37597
$ws.Cells.Item(21, 5).Value = 235275479
$ws.Cells.Item(27, 3).Value = 66243
$ws.Cells.Item(27, 4).Value = 17291
$ws.Cells.Item(27, 5).Value = 260348290
$ws.Cells.Item(30, 3).Value = 25377
$ws.Cells.Item(30, 5).Value = 106079765
$ws.Cells.Item(33, 3).Value = 3172
$ws.Cells.Item(33, 4).Value = 1218
$ws.Cells.Item(33, 5).Value = 10123368
$ws.Cells.Item(36, 3).Value = 168204
$ws.Cells.Item(36, 4).Value = 47677
$ws.Cells.Item(36, 5).Value = 297651626
$ws.Cells.Item(41, 3).Value = 93223
$ws.Cells.Item(41, 5).Value = 390579063
$ws.Cells.Item(44, 3).Value = 20970
$ws.Cells.Item(44, 4).Value = 5350
$ws.Cells.Item(44, 5).Value = 105275233
$ws.Cells.Item(47, 3).Value = 6099
$ws.Cells.Item(47, 5).Value = 22306513
$ws.Cells.Item(48, 3).Value = 120695
$ws.Cells.Item(48, 4).Value = 33166
$ws.Cells.Item(48, 5).Value = 206527569
$ws.Cells.Item(52, 3).Value = 1154
$ws.Cells.Item(52, 5).Value = 5969319
$ws.Cells.Item(54, 3).Value = 56011
$ws.Cells.Item(54, 5).Value = 213436300
$ws.Cells.Item(57, 3).Value = 23210
$ws.Cells.Item(57, 5).Value = 89989723
$ws.Cells.Item(58, 3).Value = 2566
$ws.Cells.Item(58, 4).Value = 868
$ws.Cells.Item(58, 5).Value = 7697218
$ws.Cells.Item(62, 3).Value = 38015
$ws.Cells.Item(62, 5).Value = 72534624
$ws.Cells.Item(66, 3).Value = 18896
$ws.Cells.Item(66, 5).Value = 94180400
$ws.Cells.Item(68, 3).Value = 12401
$ws.Cells.Item(68, 4).Value = 2812
$ws.Cells.Item(68, 5).Value = 56382895
$ws.Cells.Item(71, 3).Value = 260124
$ws.Cells.Item(71, 4).Value = 70788
$ws.Cells.Item(71, 5).Value = 464174291
$ws.Cells.Item(77, 3).Value = 130788
$ws.Cells.Item(77, 4).Value = 34229
$ws.Cells.Item(77, 5).Value = 517271081
$ws.Cells.Item(80, 3).Value = 63489
$ws.Cells.Item(80, 5).Value = 265752078
$ws.Cells.Item(82, 3).Value = 307
$ws.Cells.Item(82, 5).Value = 4148011
$ws.Cells.Item(83, 3).Value = 16809
$ws.Cells.Item(83, 4).Value = 6496
$ws.Cells.Item(83, 5).Value = 85876292
$ws.Cells.Item(85, 3).Value = 7435
$ws.Cells.Item(85, 4).Value = 1675
$ws.Cells.Item(85, 5).Value = 26089871
$ws.Cells.Item(87, 3).Value = 52057
$ws.Cells.Item(87, 4).Value = 12080
$ws.Cells.Item(87, 5).Value = 80862137
$ws.Cells.Item(90, 3).Value = 12269
$ws.Cells.Item(90, 5).Value = 27181742
$ws.Cells.Item(92, 3).Value = 11504
$ws.Cells.Item(92, 5).Value = 22705856
$ws.Cells.Item(94, 3).Value = 726
$ws.Cells.Item(94, 5).Value = 1253734
$ws.Cells.Item(95, 3).Value = 21341
$ws.Cells.Item(95, 5).Value = 43414595
$ws.Cells.Item(96, 3).Value = 4691
$ws.Cells.Item(96, 5).Value = 11543334
$ws.Cells.Item(98, 3).Value = 6835
$ws.Cells.Item(98, 5).Value = 16193516
$ws.Cells.Item(102, 3).Value = 252101
$ws.Cells.Item(102, 4).Value = 67182
$ws.Cells.Item(102, 5).Value = 417564322
$ws.Cells.Item(107, 3).Value = 3067
$ws.Cells.Item(107, 5).Value = 15832638
$ws.Cells.Item(109, 3).Value = 105785
$ws.Cells.Item(109, 4).Value = 27390
$ws.Cells.Item(109, 5).Value = 407605153
$ws.Cells.Item(110, 3).Value = 354
$ws.Cells.Item(110, 5).Value = 6236109
$ws.Cells.Item(112, 3).Value = 58462
$ws.Cells.Item(112, 5).Value = 232551880
$ws.Cells.Item(114, 3).Value = 36
$ws.Cells.Item(114, 5).Value = 447141
$ws.Cells.Item(115, 3).Value = 2828
$ws.Cells.Item(115, 5).Value = 8033052
$ws.Cells.Item(116, 3).Value = 5983
$ws.Cells.Item(116, 5).Value = 19448992
$ws.Cells.Item(118, 3).Value = 1015032
$ws.Cells.Item(118, 4).Value = 217561
$ws.Cells.Item(118, 5).Value = 1738604226
$ws.Cells.Item(123, 3).Value = 5385
$ws.Cells.Item(123, 5).Value = 48817590
$ws.Cells.Item(125, 3).Value = 444823
$ws.Cells.Item(125, 4).Value = 102829
$ws.Cells.Item(125, 5).Value = 1765449098
$ws.Cells.Item(128, 3).Value = 411138
$ws.Cells.Item(128, 4).Value = 87849
$ws.Cells.Item(128, 5).Value = 1615863127
$ws.Cells.Item(130, 3).Value = 5227
$ws.Cells.Item(130, 4).Value = 1613
$ws.Cells.Item(130, 5).Value = 10532108
$ws.Cells.Item(132, 3).Value = 17270
$ws.Cells.Item(132, 4).Value = 4002
$ws.Cells.Item(132, 5).Value = 60716031
$ws.Cells.Item(135, 3).Value = 62050
$ws.Cells.Item(135, 4).Value = 17299
$ws.Cells.Item(135, 5).Value = 89956936
$ws.Cells.Item(140, 3).Value = 18302
$ws.Cells.Item(140, 5).Value = 39474412
$ws.Cells.Item(147, 3).Value = 28595
$ws.Cells.Item(147, 5).Value = 44080386
$ws.Cells.Item(150, 3).Value = 11852
$ws.Cells.Item(150, 5).Value = 30350945
$ws.Cells.Item(152, 3).Value = 8460
$ws.Cells.Item(152, 4).Value = 2053
$ws.Cells.Item(152, 5).Value = 19143606
$ws.Cells.Item(155, 3).Value = 39091
$ws.Cells.Item(155, 5).Value = 94190485
$ws.Cells.Item(156, 3).Value = 3705
$ws.Cells.Item(156, 5).Value = 9782017
$ws.Cells.Item(160, 3).Value = 154782
$ws.Cells.Item(160, 4).Value = 42132
$ws.Cells.Item(160, 5).Value = 271231492
$ws.Cells.Item(167, 3).Value = 70085
$ws.Cells.Item(167, 5).Value = 286068240
$ws.Cells.Item(169, 3).Value = 28803
$ws.Cells.Item(169, 4).Value = 6990
$ws.Cells.Item(169, 5).Value = 127733650
$ws.Cells.Item(172, 3).Value = 4679
$ws.Cells.Item(172, 5).Value = 16043211
$ws.Cells.Item(173, 3).Value = 411791
$ws.Cells.Item(173, 4).Value = 114008
$ws.Cells.Item(173, 5).Value = 669491527
$ws.Cells.Item(181, 3).Value = 174251
$ws.Cells.Item(181, 5).Value = 691916500
$ws.Cells.Item(184, 3).Value = 71931
$ws.Cells.Item(184, 4).Value = 17797
$ws.Cells.Item(184, 5).Value = 302060757
$ws.Cells.Item(187, 3).Value = 10609
$ws.Cells.Item(187, 4).Value = 3867
$ws.Cells.Item(187, 5).Value = 37782830
$ws.Cells.Item(189, 3).Value = 12391
$ws.Cells.Item(189, 5).Value = 38064310
$ws.Cells.Item(191, 3).Value = 473649
$ws.Cells.Item(191, 4).Value = 125579
$ws.Cells.Item(191, 5).Value = 744969788
$ws.Cells.Item(192, 3).Value = 259
$ws.Cells.Item(192, 5).Value = 510848
$ws.Cells.Item(197, 3).Value = 1661
$ws.Cells.Item(197, 5).Value = 9619687
$ws.Cells.Item(199, 3).Value = 196623
$ws.Cells.Item(199, 4).Value = 49068
$ws.Cells.Item(199, 5).Value = 756076114
$ws.Cells.Item(202, 3).Value = 112435
$ws.Cells.Item(202, 4).Value = 26481
$ws.Cells.Item(202, 5).Value = 431661899
$ws.Cells.Item(205, 3).Value = 8742
$ws.Cells.Item(205, 4).Value = 3051
$ws.Cells.Item(205, 5).Value = 22844590
$ws.Cells.Item(208, 3).Value = 15822
$ws.Cells.Item(208, 5).Value = 47149235
$ws.Cells.Item(210, 3).Value = 183903
$ws.Cells.Item(210, 4).Value = 51825
$ws.Cells.Item(210, 5).Value = 307591593
$ws.Cells.Item(214, 3).Value = 2029
$ws.Cells.Item(214, 5).Value = 11609267
$ws.Cells.Item(216, 3).Value = 99772
$ws.Cells.Item(216, 4).Value = 26806
$ws.Cells.Item(216, 5).Value = 397667713
$ws.Cells.Item(219, 3).Value = 27988
$ws.Cells.Item(219, 5).Value = 129367685
$ws.Cells.Item(221, 3).Value = 3322
$ws.Cells.Item(221, 5).Value = 9708020
$ws.Cells.Item(223, 3).Value = 6399
$ws.Cells.Item(223, 5).Value = 22176513
$ws.Cells.Item(224, 3).Value = 482839
$ws.Cells.Item(224, 4).Value = 119801
$ws.Cells.Item(224, 5).Value = 783163787
$ws.Cells.Item(230, 3).Value = 2822
$ws.Cells.Item(230, 4).Value = 579
$ws.Cells.Item(230, 5).Value = 18119003
$ws.Cells.Item(232, 3).Value = 207248
$ws.Cells.Item(232, 4).Value = 49045
$ws.Cells.Item(232, 5).Value = 834570565
$ws.Cells.Item(235, 3).Value = 149396
$ws.Cells.Item(235, 4).Value = 33190
$ws.Cells.Item(235, 5).Value = 593970456
$ws.Cells.Item(238, 3).Value = 5136
$ws.Cells.Item(238, 5).Value = 15183092
$ws.Cells.Item(241, 3).Value = 12676
$ws.Cells.Item(241, 4).Value = 2624
$ws.Cells.Item(241, 5).Value = 39431712
